$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 305
$ws.Cells.Item(305, 4).Value = 45258
$ws.Cells.Item(305, 12).Value = 'Primera'
$ws.Cells.Item(305, 13).Value = 250
$ws.Cells.Item(305, 14).Value = 14000
$ws.Cells.Item(305, 15).Value = 14000
$ws.Cells.Item(305, 16).Value = 14000
$ws.Cells.Item(305, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(305, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(305, 19).Value = 2000

# Row 306
$ws.Cells.Item(306, 4).Value = 44475
$ws.Cells.Item(306, 12).Value = 'Primera'
$ws.Cells.Item(306, 13).Value = 200
$ws.Cells.Item(306, 14).Value = 15000
$ws.Cells.Item(306, 15).Value = 15000
$ws.Cells.Item(306, 16).Value = 15000
$ws.Cells.Item(306, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(306, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(306, 19).Value = 2143

# Row 307
$ws.Cells.Item(307, 4).Value = 45034
$ws.Cells.Item(307, 12).Value = 'Primera'
$ws.Cells.Item(307, 13).Value = 40
$ws.Cells.Item(307, 14).Value = 8000
$ws.Cells.Item(307, 15).Value = 9000
$ws.Cells.Item(307, 16).Value = 8500
$ws.Cells.Item(307, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(307, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(307, 19).Value = 1214

# Row 308
$ws.Cells.Item(308, 4).Value = 44848
$ws.Cells.Item(308, 12).Value = 'Primera'
$ws.Cells.Item(308, 13).Value = 450
$ws.Cells.Item(308, 14).Value = 10000
$ws.Cells.Item(308, 15).Value = 13000
$ws.Cells.Item(308, 16).Value = 11333
$ws.Cells.Item(308, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(308, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(308, 19).Value = 1619

# Row 309
$ws.Cells.Item(309, 4).Value = 44623
$ws.Cells.Item(309, 12).Value = 'Primera'
$ws.Cells.Item(309, 13).Value = 105
$ws.Cells.Item(309, 14).Value = 6000
$ws.Cells.Item(309, 15).Value = 7000
$ws.Cells.Item(309, 16).Value = 6619
$ws.Cells.Item(309, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(309, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(309, 19).Value = 946

# Row 310
$ws.Cells.Item(310, 4).Value = 44284
$ws.Cells.Item(310, 12).Value = 'Primera'
$ws.Cells.Item(310, 13).Value = 65
$ws.Cells.Item(310, 14).Value = 7000
$ws.Cells.Item(310, 15).Value = 7000
$ws.Cells.Item(310, 16).Value = 7000
$ws.Cells.Item(310, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(310, 18).Value = 'Provincia de Cautín'
$ws.Cells.Item(310, 19).Value = 1000

# Row 311
$ws.Cells.Item(311, 4).Value = 44274
$ws.Cells.Item(311, 12).Value = 'Primera'
$ws.Cells.Item(311, 13).Value = 95
$ws.Cells.Item(311, 14).Value = 6000
$ws.Cells.Item(311, 15).Value = 7000
$ws.Cells.Item(311, 16).Value = 6579
$ws.Cells.Item(311, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(311, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(311, 19).Value = 940

# Row 312
$ws.Cells.Item(312, 4).Value = 45243
$ws.Cells.Item(312, 12).Value = 'Primera'
$ws.Cells.Item(312, 13).Value = 1200
$ws.Cells.Item(312, 14).Value = 15000
$ws.Cells.Item(312, 15).Value = 16000
$ws.Cells.Item(312, 16).Value = 15333
$ws.Cells.Item(312, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(312, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(312, 19).Value = 2190

# Row 313
$ws.Cells.Item(313, 4).Value = 45243
$ws.Cells.Item(313, 12).Value = 'Segunda'
$ws.Cells.Item(313, 13).Value = 100
$ws.Cells.Item(313, 14).Value = 14000
$ws.Cells.Item(313, 15).Value = 14000
$ws.Cells.Item(313, 16).Value = 14000
$ws.Cells.Item(313, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(313, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(313, 19).Value = 2000

# Row 314
$ws.Cells.Item(314, 4).Value = 44858
$ws.Cells.Item(314, 12).Value = 'Especial'
$ws.Cells.Item(314, 13).Value = 400
$ws.Cells.Item(314, 14).Value = 13000
$ws.Cells.Item(314, 15).Value = 13000
$ws.Cells.Item(314, 16).Value = 13000
$ws.Cells.Item(314, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(314, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(314, 19).Value = 1857

# Row 315
$ws.Cells.Item(315, 4).Value = 44858
$ws.Cells.Item(315, 12).Value = 'Primera'
$ws.Cells.Item(315, 13).Value = 1400
$ws.Cells.Item(315, 14).Value = 9000
$ws.Cells.Item(315, 15).Value = 10000
$ws.Cells.Item(315, 16).Value = 9571
$ws.Cells.Item(315, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(315, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(315, 19).Value = 1367

# Row 316
$ws.Cells.Item(316, 4).Value = 44858
$ws.Cells.Item(316, 12).Value = 'Segunda'
$ws.Cells.Item(316, 13).Value = 200
$ws.Cells.Item(316, 14).Value = 7000
$ws.Cells.Item(316, 15).Value = 7000
$ws.Cells.Item(316, 16).Value = 7000
$ws.Cells.Item(316, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(316, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(316, 19).Value = 1000

# Row 317
$ws.Cells.Item(317, 4).Value = 44484
$ws.Cells.Item(317, 12).Value = 'Primera'
$ws.Cells.Item(317, 13).Value = 220
$ws.Cells.Item(317, 14).Value = 11000
$ws.Cells.Item(317, 15).Value = 12000
$ws.Cells.Item(317, 16).Value = 11432
$ws.Cells.Item(317, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(317, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(317, 19).Value = 1633

# Row 318
$ws.Cells.Item(318, 4).Value = 44484
$ws.Cells.Item(318, 12).Value = 'Segunda'
$ws.Cells.Item(318, 13).Value = 125
$ws.Cells.Item(318, 14).Value = 7000
$ws.Cells.Item(318, 15).Value = 7000
$ws.Cells.Item(318, 16).Value = 7000
$ws.Cells.Item(318, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(318, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(318, 19).Value = 1000

# Row 319
$ws.Cells.Item(319, 4).Value = 44876
$ws.Cells.Item(319, 12).Value = 'Primera'
$ws.Cells.Item(319, 13).Value = 330
$ws.Cells.Item(319, 14).Value = 7000
$ws.Cells.Item(319, 15).Value = 8000
$ws.Cells.Item(319, 16).Value = 7455
$ws.Cells.Item(319, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(319, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(319, 19).Value = 1065

# Row 320
$ws.Cells.Item(320, 4).Value = 44876
$ws.Cells.Item(320, 12).Value = 'Primera'
$ws.Cells.Item(320, 13).Value = 100
$ws.Cells.Item(320, 14).Value = 8000
$ws.Cells.Item(320, 15).Value = 8000
$ws.Cells.Item(320, 16).Value = 8000
$ws.Cells.Item(320, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(320, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(320, 19).Value = 1143

# Row 321
$ws.Cells.Item(321, 4).Value = 44876
$ws.Cells.Item(321, 12).Value = 'Segunda'
$ws.Cells.Item(321, 13).Value = 50
$ws.Cells.Item(321, 14).Value = 6000
$ws.Cells.Item(321, 15).Value = 6000
$ws.Cells.Item(321, 16).Value = 6000
$ws.Cells.Item(321, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(321, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(321, 19).Value = 857

# Row 322
$ws.Cells.Item(322, 4).Value = 44522
$ws.Cells.Item(322, 12).Value = 'Primera'
$ws.Cells.Item(322, 13).Value = 3300
$ws.Cells.Item(322, 14).Value = 9000
$ws.Cells.Item(322, 15).Value = 10000
$ws.Cells.Item(322, 16).Value = 9455
$ws.Cells.Item(322, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(322, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(322, 19).Value = 1351

# Row 323
$ws.Cells.Item(323, 4).Value = 44522
$ws.Cells.Item(323, 12).Value = 'Primera'
$ws.Cells.Item(323, 13).Value = 400
$ws.Cells.Item(323, 14).Value = 8000
$ws.Cells.Item(323, 15).Value = 8000
$ws.Cells.Item(323, 16).Value = 8000
$ws.Cells.Item(323, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(323, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(323, 19).Value = 1143

# Row 324
$ws.Cells.Item(324, 4).Value = 44581
$ws.Cells.Item(324, 12).Value = 'Especial'
$ws.Cells.Item(324, 13).Value = 250
$ws.Cells.Item(324, 14).Value = 10000
$ws.Cells.Item(324, 15).Value = 10000
$ws.Cells.Item(324, 16).Value = 10000
$ws.Cells.Item(324, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(324, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(324, 19).Value = 1429

# Row 325
$ws.Cells.Item(325, 4).Value = 44581
$ws.Cells.Item(325, 12).Value = 'Primera'
$ws.Cells.Item(325, 13).Value = 185
$ws.Cells.Item(325, 14).Value = 7000
$ws.Cells.Item(325, 15).Value = 7000
$ws.Cells.Item(325, 16).Value = 7000
$ws.Cells.Item(325, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(325, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(325, 19).Value = 1000

# Row 326
$ws.Cells.Item(326, 4).Value = 44581
$ws.Cells.Item(326, 12).Value = 'Primera'
$ws.Cells.Item(326, 13).Value = 85
$ws.Cells.Item(326, 14).Value = 7000
$ws.Cells.Item(326, 15).Value = 7000
$ws.Cells.Item(326, 16).Value = 7000
$ws.Cells.Item(326, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(326, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(326, 19).Value = 1000

# Row 327
$ws.Cells.Item(327, 4).Value = 44581
$ws.Cells.Item(327, 12).Value = 'Segunda'
$ws.Cells.Item(327, 13).Value = 110
$ws.Cells.Item(327, 14).Value = 5000
$ws.Cells.Item(327, 15).Value = 5000
$ws.Cells.Item(327, 16).Value = 5000
$ws.Cells.Item(327, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(327, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(327, 19).Value = 714

# Row 328
$ws.Cells.Item(328, 4).Value = 44532
$ws.Cells.Item(328, 12).Value = 'Primera'
$ws.Cells.Item(328, 13).Value = 2000
$ws.Cells.Item(328, 14).Value = 8000
$ws.Cells.Item(328, 15).Value = 8000
$ws.Cells.Item(328, 16).Value = 8000
$ws.Cells.Item(328, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(328, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(328, 19).Value = 1143

# Row 329
$ws.Cells.Item(329, 4).Value = 44532
$ws.Cells.Item(329, 12).Value = 'Primera'
$ws.Cells.Item(329, 13).Value = 500
$ws.Cells.Item(329, 14).Value = 8000
$ws.Cells.Item(329, 15).Value = 8000
$ws.Cells.Item(329, 16).Value = 8000
$ws.Cells.Item(329, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(329, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(329, 19).Value = 1143

# Row 330
$ws.Cells.Item(330, 4).Value = 44258
$ws.Cells.Item(330, 12).Value = 'Primera'
$ws.Cells.Item(330, 13).Value = 65
$ws.Cells.Item(330, 14).Value = 7000
$ws.Cells.Item(330, 15).Value = 7000
$ws.Cells.Item(330, 16).Value = 7000
$ws.Cells.Item(330, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(330, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(330, 19).Value = 1000

# Row 331
$ws.Cells.Item(331, 4).Value = 44454
$ws.Cells.Item(331, 12).Value = 'Tercera'
$ws.Cells.Item(331, 13).Value = 80
$ws.Cells.Item(331, 14).Value = 18000
$ws.Cells.Item(331, 15).Value = 18000
$ws.Cells.Item(331, 16).Value = 18000
$ws.Cells.Item(331, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(331, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(331, 19).Value = 2571

# Row 332
$ws.Cells.Item(332, 4).Value = 44232
$ws.Cells.Item(332, 12).Value = 'Primera'
$ws.Cells.Item(332, 13).Value = 100
$ws.Cells.Item(332, 14).Value = 7000
$ws.Cells.Item(332, 15).Value = 8000
$ws.Cells.Item(332, 16).Value = 7550
$ws.Cells.Item(332, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(332, 18).Value = 'Provincia de Cautín'
$ws.Cells.Item(332, 19).Value = 1079

# Row 333
$ws.Cells.Item(333, 4).Value = 44832
$ws.Cells.Item(333, 12).Value = 'Primera'
$ws.Cells.Item(333, 13).Value = 400
$ws.Cells.Item(333, 14).Value = 17000
$ws.Cells.Item(333, 15).Value = 17000
$ws.Cells.Item(333, 16).Value = 17000
$ws.Cells.Item(333, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(333, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(333, 19).Value = 2429

# Row 334
$ws.Cells.Item(334, 4).Value = 44364
$ws.Cells.Item(334, 12).Value = 'Primera'
$ws.Cells.Item(334, 13).Value = 40
$ws.Cells.Item(334, 14).Value = 20000
$ws.Cells.Item(334, 15).Value = 20000
$ws.Cells.Item(334, 16).Value = 20000
$ws.Cells.Item(334, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(334, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(334, 19).Value = 2857

# Row 335
$ws.Cells.Item(335, 4).Value = 44600
$ws.Cells.Item(335, 12).Value = 'Primera'
$ws.Cells.Item(335, 13).Value = 80
$ws.Cells.Item(335, 14).Value = 8000
$ws.Cells.Item(335, 15).Value = 8000
$ws.Cells.Item(335, 16).Value = 8000
$ws.Cells.Item(335, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(335, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(335, 19).Value = 1143

# Row 336
$ws.Cells.Item(336, 4).Value = 44246
$ws.Cells.Item(336, 12).Value = 'Primera'
$ws.Cells.Item(336, 13).Value = 60
$ws.Cells.Item(336, 14).Value = 7000
$ws.Cells.Item(336, 15).Value = 8000
$ws.Cells.Item(336, 16).Value = 7417
$ws.Cells.Item(336, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(336, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(336, 19).Value = 1060

# Row 337
$ws.Cells.Item(337, 4).Value = 44246
$ws.Cells.Item(337, 12).Value = 'Segunda'
$ws.Cells.Item(337, 13).Value = 50
$ws.Cells.Item(337, 14).Value = 5000
$ws.Cells.Item(337, 15).Value = 5000
$ws.Cells.Item(337, 16).Value = 5000
$ws.Cells.Item(337, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(337, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(337, 19).Value = 714

# Row 338
$ws.Cells.Item(338, 4).Value = 44491
$ws.Cells.Item(338, 12).Value = 'Primera'
$ws.Cells.Item(338, 13).Value = 500
$ws.Cells.Item(338, 14).Value = 8000
$ws.Cells.Item(338, 15).Value = 9000
$ws.Cells.Item(338, 16).Value = 8600
$ws.Cells.Item(338, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(338, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(338, 19).Value = 1229

# Row 339
$ws.Cells.Item(339, 4).Value = 44491
$ws.Cells.Item(339, 12).Value = 'Segunda'
$ws.Cells.Item(339, 13).Value = 100
$ws.Cells.Item(339, 14).Value = 6000
$ws.Cells.Item(339, 15).Value = 6000
$ws.Cells.Item(339, 16).Value = 6000
$ws.Cells.Item(339, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(339, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(339, 19).Value = 857

# Row 340
$ws.Cells.Item(340, 4).Value = 44179
$ws.Cells.Item(340, 12).Value = 'Primera'
$ws.Cells.Item(340, 13).Value = 600
$ws.Cells.Item(340, 14).Value = 10000
$ws.Cells.Item(340, 15).Value = 11000
$ws.Cells.Item(340, 16).Value = 10500
$ws.Cells.Item(340, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(340, 18).Value = 'Región del Maule'
$ws.Cells.Item(340, 19).Value = 1500

# Row 341
$ws.Cells.Item(341, 4).Value = 44179
$ws.Cells.Item(341, 12).Value = 'Primera'
$ws.Cells.Item(341, 13).Value = 200
$ws.Cells.Item(341, 14).Value = 10000
$ws.Cells.Item(341, 15).Value = 10000
$ws.Cells.Item(341, 16).Value = 10000
$ws.Cells.Item(341, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(341, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(341, 19).Value = 1429

# Row 342
$ws.Cells.Item(342, 4).Value = 44181
$ws.Cells.Item(342, 12).Value = 'Primera'
$ws.Cells.Item(342, 13).Value = 260
$ws.Cells.Item(342, 14).Value = 9000
$ws.Cells.Item(342, 15).Value = 10000
$ws.Cells.Item(342, 16).Value = 9462
$ws.Cells.Item(342, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(342, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(342, 19).Value = 1352

# Row 343
$ws.Cells.Item(343, 4).Value = 44312
$ws.Cells.Item(343, 12).Value = 'Primera'
$ws.Cells.Item(343, 13).Value = 65
$ws.Cells.Item(343, 14).Value = 7000
$ws.Cells.Item(343, 15).Value = 7000
$ws.Cells.Item(343, 16).Value = 7000
$ws.Cells.Item(343, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(343, 18).Value = 'Provincia de Cautín'
$ws.Cells.Item(343, 19).Value = 1000

# Row 344
$ws.Cells.Item(344, 4).Value = 44270
$ws.Cells.Item(344, 12).Value = 'Primera'
$ws.Cells.Item(344, 13).Value = 65
$ws.Cells.Item(344, 14).Value = 7000
$ws.Cells.Item(344, 15).Value = 7000
$ws.Cells.Item(344, 16).Value = 7000
$ws.Cells.Item(344, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(344, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(344, 19).Value = 1000

# Row 345
$ws.Cells.Item(345, 4).Value = 44277
$ws.Cells.Item(345, 12).Value = 'Primera'
$ws.Cells.Item(345, 13).Value = 80
$ws.Cells.Item(345, 14).Value = 6000
$ws.Cells.Item(345, 15).Value = 6000
$ws.Cells.Item(345, 16).Value = 6000
$ws.Cells.Item(345, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(345, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(345, 19).Value = 857

# Row 346
$ws.Cells.Item(346, 4).Value = 45222
$ws.Cells.Item(346, 12).Value = 'Primera'
$ws.Cells.Item(346, 13).Value = 95
$ws.Cells.Item(346, 14).Value = 10000
$ws.Cells.Item(346, 15).Value = 10000
$ws.Cells.Item(346, 16).Value = 10000
$ws.Cells.Item(346, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(346, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(346, 19).Value = 1429

# Row 347
$ws.Cells.Item(347, 4).Value = 44592
$ws.Cells.Item(347, 12).Value = 'Primera'
$ws.Cells.Item(347, 13).Value = 55
$ws.Cells.Item(347, 14).Value = 8000
$ws.Cells.Item(347, 15).Value = 8000
$ws.Cells.Item(347, 16).Value = 8000
$ws.Cells.Item(347, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(347, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(347, 19).Value = 1143

# Row 348
$ws.Cells.Item(348, 4).Value = 44714
$ws.Cells.Item(348, 12).Value = 'Especial'
$ws.Cells.Item(348, 13).Value = 600
$ws.Cells.Item(348, 14).Value = 18000
$ws.Cells.Item(348, 15).Value = 18000
$ws.Cells.Item(348, 16).Value = 18000
$ws.Cells.Item(348, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(348, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(348, 19).Value = 2571

# Row 349
$ws.Cells.Item(349, 4).Value = 44271
$ws.Cells.Item(349, 12).Value = 'Primera'
$ws.Cells.Item(349, 13).Value = 65
$ws.Cells.Item(349, 14).Value = 7000
$ws.Cells.Item(349, 15).Value = 7000
$ws.Cells.Item(349, 16).Value = 7000
$ws.Cells.Item(349, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(349, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(349, 19).Value = 1000

# Row 350
$ws.Cells.Item(350, 4).Value = 44880
$ws.Cells.Item(350, 12).Value = 'Primera'
$ws.Cells.Item(350, 13).Value = 125
$ws.Cells.Item(350, 14).Value = 9000
$ws.Cells.Item(350, 15).Value = 9000
$ws.Cells.Item(350, 16).Value = 9000
$ws.Cells.Item(350, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(350, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(350, 19).Value = 1286

# Row 351
$ws.Cells.Item(351, 4).Value = 44880
$ws.Cells.Item(351, 12).Value = 'Segunda'
$ws.Cells.Item(351, 13).Value = 95
$ws.Cells.Item(351, 14).Value = 6500
$ws.Cells.Item(351, 15).Value = 6500
$ws.Cells.Item(351, 16).Value = 6500
$ws.Cells.Item(351, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(351, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(351, 19).Value = 929

# Row 352
$ws.Cells.Item(352, 4).Value = 44399
$ws.Cells.Item(352, 12).Value = 'Primera'
$ws.Cells.Item(352, 13).Value = 65
$ws.Cells.Item(352, 14).Value = 22000
$ws.Cells.Item(352, 15).Value = 23000
$ws.Cells.Item(352, 16).Value = 22538
$ws.Cells.Item(352, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(352, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(352, 19).Value = 3220

# Row 353
$ws.Cells.Item(353, 4).Value = 44161
$ws.Cells.Item(353, 12).Value = 'Primera'
$ws.Cells.Item(353, 13).Value = 550
$ws.Cells.Item(353, 14).Value = 9000
$ws.Cells.Item(353, 15).Value = 10000
$ws.Cells.Item(353, 16).Value = 9455
$ws.Cells.Item(353, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(353, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(353, 19).Value = 1351

# Row 354
$ws.Cells.Item(354, 4).Value = 44161
$ws.Cells.Item(354, 12).Value = 'Primera'
$ws.Cells.Item(354, 13).Value = 305
$ws.Cells.Item(354, 14).Value = 6500
$ws.Cells.Item(354, 15).Value = 7000
$ws.Cells.Item(354, 16).Value = 6705
$ws.Cells.Item(354, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(354, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(354, 19).Value = 958

# Row 355
$ws.Cells.Item(355, 4).Value = 44161
$ws.Cells.Item(355, 12).Value = 'Segunda'
$ws.Cells.Item(355, 13).Value = 180
$ws.Cells.Item(355, 14).Value = 4000
$ws.Cells.Item(355, 15).Value = 4000
$ws.Cells.Item(355, 16).Value = 4000
$ws.Cells.Item(355, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(355, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(355, 19).Value = 571

# Row 356
$ws.Cells.Item(356, 4).Value = 44509
$ws.Cells.Item(356, 12).Value = 'Primera'
$ws.Cells.Item(356, 13).Value = 270
$ws.Cells.Item(356, 14).Value = 9000
$ws.Cells.Item(356, 15).Value = 9000
$ws.Cells.Item(356, 16).Value = 9000
$ws.Cells.Item(356, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(356, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(356, 19).Value = 1286

# Row 357
$ws.Cells.Item(357, 4).Value = 44547
$ws.Cells.Item(357, 12).Value = 'Primera'
$ws.Cells.Item(357, 13).Value = 325
$ws.Cells.Item(357, 14).Value = 8000
$ws.Cells.Item(357, 15).Value = 9000
$ws.Cells.Item(357, 16).Value = 8615
$ws.Cells.Item(357, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(357, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(357, 19).Value = 1231

# Row 358
$ws.Cells.Item(358, 4).Value = 44547
$ws.Cells.Item(358, 12).Value = 'Segunda'
$ws.Cells.Item(358, 13).Value = 185
$ws.Cells.Item(358, 14).Value = 7000
$ws.Cells.Item(358, 15).Value = 7000
$ws.Cells.Item(358, 16).Value = 7000
$ws.Cells.Item(358, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(358, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(358, 19).Value = 1000

# Row 359
$ws.Cells.Item(359, 4).Value = 44539
$ws.Cells.Item(359, 12).Value = 'Primera'
$ws.Cells.Item(359, 13).Value = 350
$ws.Cells.Item(359, 14).Value = 9000
$ws.Cells.Item(359, 15).Value = 10000
$ws.Cells.Item(359, 16).Value = 9557
$ws.Cells.Item(359, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(359, 18).Value = 'Región del Maule'
$ws.Cells.Item(359, 19).Value = 1365

# Row 360
$ws.Cells.Item(360, 4).Value = 44539
$ws.Cells.Item(360, 12).Value = 'Primera'
$ws.Cells.Item(360, 13).Value = 90
$ws.Cells.Item(360, 14).Value = 8000
$ws.Cells.Item(360, 15).Value = 9000
$ws.Cells.Item(360, 16).Value = 8389
$ws.Cells.Item(360, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(360, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(360, 19).Value = 1198

# Row 361
$ws.Cells.Item(361, 4).Value = 44193
$ws.Cells.Item(361, 12).Value = 'Primera'
$ws.Cells.Item(361, 13).Value = 150
$ws.Cells.Item(361, 14).Value = 10000
$ws.Cells.Item(361, 15).Value = 10000
$ws.Cells.Item(361, 16).Value = 10000
$ws.Cells.Item(361, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(361, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(361, 19).Value = 1429

# Row 362
$ws.Cells.Item(362, 4).Value = 44193
$ws.Cells.Item(362, 12).Value = 'Segunda'
$ws.Cells.Item(362, 13).Value = 90
$ws.Cells.Item(362, 14).Value = 9000
$ws.Cells.Item(362, 15).Value = 9000
$ws.Cells.Item(362, 16).Value = 9000
$ws.Cells.Item(362, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(362, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(362, 19).Value = 1286

# Row 363
$ws.Cells.Item(363, 4).Value = 44917
$ws.Cells.Item(363, 12).Value = 'Primera'
$ws.Cells.Item(363, 13).Value = 300
$ws.Cells.Item(363, 14).Value = 8500
$ws.Cells.Item(363, 15).Value = 8500
$ws.Cells.Item(363, 16).Value = 8500
$ws.Cells.Item(363, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(363, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(363, 19).Value = 1214

# Row 364
$ws.Cells.Item(364, 4).Value = 44914
$ws.Cells.Item(364, 12).Value = 'Primera'
$ws.Cells.Item(364, 13).Value = 45
$ws.Cells.Item(364, 14).Value = 9000
$ws.Cells.Item(364, 15).Value = 9000
$ws.Cells.Item(364, 16).Value = 9000
$ws.Cells.Item(364, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(364, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(364, 19).Value = 1286

# Row 365
$ws.Cells.Item(365, 4).Value = 44574
$ws.Cells.Item(365, 12).Value = 'Primera'
$ws.Cells.Item(365, 13).Value = 100
$ws.Cells.Item(365, 14).Value = 8000
$ws.Cells.Item(365, 15).Value = 8000
$ws.Cells.Item(365, 16).Value = 8000
$ws.Cells.Item(365, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(365, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(365, 19).Value = 1143

# Row 366
$ws.Cells.Item(366, 4).Value = 44214
$ws.Cells.Item(366, 12).Value = 'Primera'
$ws.Cells.Item(366, 13).Value = 55
$ws.Cells.Item(366, 14).Value = 7000
$ws.Cells.Item(366, 15).Value = 8000
$ws.Cells.Item(366, 16).Value = 7455
$ws.Cells.Item(366, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(366, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(366, 19).Value = 1065

# Row 367
$ws.Cells.Item(367, 4).Value = 44214
$ws.Cells.Item(367, 12).Value = 'Segunda'
$ws.Cells.Item(367, 13).Value = 40
$ws.Cells.Item(367, 14).Value = 4500
$ws.Cells.Item(367, 15).Value = 4500
$ws.Cells.Item(367, 16).Value = 4500
$ws.Cells.Item(367, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(367, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(367, 19).Value = 643

# Row 368
$ws.Cells.Item(368, 4).Value = 44567
$ws.Cells.Item(368, 12).Value = 'Primera'
$ws.Cells.Item(368, 13).Value = 130
$ws.Cells.Item(368, 14).Value = 7000
$ws.Cells.Item(368, 15).Value = 8000
$ws.Cells.Item(368, 16).Value = 7538
$ws.Cells.Item(368, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(368, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(368, 19).Value = 1077

# Row 369
$ws.Cells.Item(369, 4).Value = 44203
$ws.Cells.Item(369, 12).Value = 'Primera'
$ws.Cells.Item(369, 13).Value = 175
$ws.Cells.Item(369, 14).Value = 7000
$ws.Cells.Item(369, 15).Value = 8000
$ws.Cells.Item(369, 16).Value = 7543
$ws.Cells.Item(369, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(369, 18).Value = 'Provincia de Cautín'
$ws.Cells.Item(369, 19).Value = 1078

# Row 370
$ws.Cells.Item(370, 4).Value = 44203
$ws.Cells.Item(370, 12).Value = 'Segunda'
$ws.Cells.Item(370, 13).Value = 65
$ws.Cells.Item(370, 14).Value = 5500
$ws.Cells.Item(370, 15).Value = 5500
$ws.Cells.Item(370, 16).Value = 5500
$ws.Cells.Item(370, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(370, 18).Value = 'Provincia de Cautín'
$ws.Cells.Item(370, 19).Value = 786

# Row 371 (new row)
$ws.Cells.Item(371, 1).Value = 10
$ws.Cells.Item(371, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(371, 3).Value = 'La Araucanía'
$ws.Cells.Item(371, 4).Value = 44189
$ws.Cells.Item(371, 5).Value = 9
$ws.Cells.Item(371, 6).Value = 'Fruta'
$ws.Cells.Item(371, 7).Value = 100101
$ws.Cells.Item(371, 8).Value = 'Berries'
$ws.Cells.Item(371, 9).Value = 100112025
$ws.Cells.Item(371, 10).Value = 'Frutilla'
$ws.Cells.Item(371, 11).Value = 'Sin especificar'
$ws.Cells.Item(371, 12).Value = 'Primera'
$ws.Cells.Item(371, 13).Value = 400
$ws.Cells.Item(371, 14).Value = 8000
$ws.Cells.Item(371, 15).Value = 10000
$ws.Cells.Item(371, 16).Value = 9000
$ws.Cells.Item(371, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(371, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(371, 19).Value = 1286
$ws.Cells.Item(371, 20).Value = 7
